# "need to work on fin analysis"
# Refresh the P&L / Trading-PL workpaper: roll the report date from
# 01-May-2024 to 01-Jan-2025, update every figure that moved, drop a
# couple of now-unused entries (B11/B18), and append the new
# "Trading PL" block (Conveyance Charges, SALARY & WAGES, Commission on
# sales, TRAVELLING CHARGES) at rows 50-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "(Trading PL)" section labels first, so the shared-string table
# --- gets the same append order as the authored workbook.
$ws.Range("D50").Value = "Conveyance Charges"
$ws.Range("D51").Value = "SALARY & WAGES"
$ws.Range("D52").Value = "Commission on sales"
$ws.Range("D53").Value = "TRAVELLING CHARGES"
$ws.Range("C50").Value = "(Trading PL)"

# --- Report date (B1/E1) ---
$ws.Range("B1").Value = 45658
$ws.Range("E1").Value = 45658

# --- Updated figures, column E (right-hand P&L block) ---
$ws.Range("E2").Value = 261259
$ws.Range("E3").Value = 607103
$ws.Range("E4").Value = 6280941.5
$ws.Range("E5").Value = 761934
$ws.Range("E6").Value = 452225.5
$ws.Range("E7").Value = 80099.679999999993
$ws.Range("E8").Value = 693944
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 37866
$ws.Range("E11").Value = 478873
$ws.Range("E13").Value = 4970
$ws.Range("E14").Value = 9659215.6799999997
$ws.Range("E15").Value = 10643.6
$ws.Range("E16").Value = 193459.81
$ws.Range("E17").Value = 524196
$ws.Range("E18").Value = 728299.41
$ws.Range("E19").Value = 10387515.09
$ws.Range("E22").Value = 608383
$ws.Range("E23").Value = 355700
$ws.Range("E26").Value = 532489.14
$ws.Range("E27").Value = 53840.53
$ws.Range("E28").Value = 0
$ws.Range("E29").Value = 216333
$ws.Range("E30").Value = 2416745.67
$ws.Range("E31").Value = 12804260.76
$ws.Range("E32").Value = 12804261
$ws.Range("E33").Value = -0.24000000022351742
$ws.Range("E34").Value = 1832445
$ws.Range("E35").Value = 3119379
$ws.Range("E37").Value = 11273298.68
$ws.Range("E38").Value = 1530962.08
$ws.Range("E39").Value = 12804260.76
$ws.Range("E40").Value = 11273298.68
$ws.Range("E41").Value = 1530962.08
$ws.Range("E42").Value = 12804260.76
$ws.Range("E44").Value = 3140779.59
$ws.Range("E45").Value = 7149303.5
$ws.Range("E46").Value = 97432
$ws.Range("E47").Value = 2291745.67
$ws.Range("E49").Value = 12804260.76
$ws.Range("E50").Value = 3432
$ws.Range("E51").Value = 94000

# --- Updated figures, column B (left-hand manufacturing block) ---
$ws.Range("B3").Value = 3119379
$ws.Range("B4").Value = 4970
$ws.Range("B6").Value = 6280941.5
$ws.Range("B8").Value = 37866
$ws.Range("B10").Value = 80099.679999999993
$ws.Range("B15").Value = 693944
$ws.Range("B17").Value = 478873
$ws.Range("B19").Value = 452225.5
$ws.Range("B21").Value = 11273298.68
$ws.Range("B23").Value = 4150113
$ws.Range("B25").Value = 717244
$ws.Range("B26").Value = 6280942
$ws.Range("B27").Value = 11148299
$ws.Range("B29").Value = 11273299
$ws.Range("B30").Value = 6998186
$ws.Range("B32").Value = 9344.7000000000007
$ws.Range("B33").Value = 74639.199999999997
$ws.Range("B34").Value = 79408
$ws.Range("B35").Value = 11352707
$ws.Range("B36").Value = 22626006
$ws.Range("B38").Value = 11273298.68

# --- Entries that are no longer populated ---
$ws.Range("B11").ClearContents()
$ws.Range("B18").ClearContents()

# --- Cursor / selection moves back to the top of the updated column ---
$ws.Range("E2").Select() | Out-Null

# --- Page is now set up for printing in portrait orientation ---
$ws.PageSetup.Orientation = 1
